$wb = $excel.ActiveWorkbook

# --- Players sheet ---
$players = $wb.Worksheets.Item("Players")

# Row 8
$players.Range("D8").Value = 'Labaron Philon Jr.'
$players.Range("E8").Value = 'ALA'
$players.Range("H8").Value = 12
$players.Range("I8").Value = 18
$players.Range("J8").Value = 1
$players.Range("K8").Value = 1
$players.Range("L8").Value = 1
$players.Range("M8").Value = 0
$players.Range("N8").Value = 4
$players.Range("O8").Value = 20
# Row 9
$players.Range("D9").Value = 'Derrion Reid'
$players.Range("E9").Value = 'OU'
$players.Range("F9").Value = 'OU@MSST'
$players.Range("H9").Value = 8
$players.Range("I9").Value = 11
$players.Range("J9").Value = 3
$players.Range("K9").Value = 0
$players.Range("L9").Value = 1
$players.Range("N9").Value = 1
$players.Range("O9").Value = 22
# Row 10
$players.Range("D10").Value = 'Devin McGlockton'
$players.Range("E10").Value = 'VAN'
$players.Range("F10").Value = 'ALA@VAN'
$players.Range("H10").Value = 21
$players.Range("I10").Value = 10
$players.Range("J10").Value = 12
$players.Range("K10").Value = 1
$players.Range("M10").Value = 1
$players.Range("N10").Value = 1
$players.Range("O10").Value = 29
# Row 11
$players.Range("D11").Value = 'Billy Richmond III'
$players.Range("E11").Value = 'ARK'
$players.Range("H11").Value = 20
$players.Range("I11").Value = 13
$players.Range("J11").Value = 6
$players.Range("K11").Value = 1
$players.Range("L11").Value = 0
$players.Range("M11").Value = 0
$players.Range("N11").Value = 0
$players.Range("O11").Value = 24
# Row 12
$players.Range("D12").Value = 'James Scott'
$players.Range("E12").Value = 'MISS'
$players.Range("F12").Value = 'ARK@MISS'
$players.Range("H12").Value = 11
$players.Range("I12").Value = 6
$players.Range("J12").Value = 3
$players.Range("K12").Value = 2
$players.Range("L12").Value = 0
$players.Range("M12").Value = 1
$players.Range("N12").Value = 1
$players.Range("O12").Value = 24
# Row 16
$players.Range("D16").Value = 'Meechie Johnson'
$players.Range("E16").Value = 'SC'
$players.Range("F16").Value = 'SC@LSU'
$players.Range("H16").Value = 22
$players.Range("I16").Value = 19
$players.Range("J16").Value = 2
$players.Range("K16").Value = 6
$players.Range("L16").Value = 3
$players.Range("M16").Value = 0
$players.Range("N16").Value = 1
$players.Range("O16").Value = 35
# Row 17
$players.Range("C17").Value = 'Yes'
$players.Range("D17").Value = 'KeShawn Murphy'
$players.Range("E17").Value = 'AUB'
$players.Range("F17").Value = 'TA&M@AUB'
$players.Range("H17").Value = 12
$players.Range("I17").Value = 8
$players.Range("J17").Value = 8
$players.Range("K17").Value = 1
$players.Range("L17").Value = 0
$players.Range("M17").Value = 1
$players.Range("N17").Value = 1
$players.Range("O17").Value = 27
# Row 21
$players.Range("C21").Value = 'No'
$players.Range("D21").Value = 'Malik Dia'
$players.Range("E21").Value = 'MISS'
$players.Range("H21").Value = 18
$players.Range("I21").Value = 16
$players.Range("J21").Value = 6
$players.Range("K21").Value = 1
$players.Range("L21").Value = 2
$players.Range("M21").Value = 2
$players.Range("N21").Value = 2
$players.Range("O21").Value = 27
# Row 22
$players.Range("D22").Value = 'Denzel Aberdeen'
$players.Range("E22").Value = 'UK'
$players.Range("F22").Value = 'MIZ@UK'
$players.Range("H22").Value = 9
$players.Range("I22").Value = 7
$players.Range("J22").Value = 3
$players.Range("K22").Value = 4
$players.Range("L22").Value = 0
$players.Range("M22").Value = 0
$players.Range("N22").Value = 1
$players.Range("O22").Value = 28
# Row 23
$players.Range("D23").Value = 'Meleek Thomas'
$players.Range("E23").Value = 'ARK'
$players.Range("F23").Value = 'ARK@MISS'
$players.Range("H23").Value = 8
$players.Range("I23").Value = 13
$players.Range("J23").Value = 1
$players.Range("K23").Value = 2
$players.Range("L23").Value = 2
$players.Range("M23").Value = 0
$players.Range("N23").Value = 2
$players.Range("O23").Value = 25
# Row 30
$players.Range("C30").Value = 'No'
$players.Range("D30").Value = 'Ilias Kamardine'
$players.Range("E30").Value = 'MISS'
$players.Range("F30").Value = 'ARK@MISS'
$players.Range("H30").Value = 17
$players.Range("I30").Value = 16
$players.Range("J30").Value = 4
$players.Range("K30").Value = 4
$players.Range("L30").Value = 3
$players.Range("M30").Value = 0
$players.Range("N30").Value = 4
$players.Range("O30").Value = 26
# Row 31
$players.Range("D31").Value = 'Quincy Ballard'
$players.Range("E31").Value = 'MSST'
$players.Range("F31").Value = 'OU@MSST'
$players.Range("H31").Value = 7
$players.Range("I31").Value = 10
$players.Range("J31").Value = 4
$players.Range("K31").Value = 0
$players.Range("L31").Value = 0
$players.Range("M31").Value = 1
$players.Range("N31").Value = 2
$players.Range("O31").Value = 16
# Row 60
$players.Range("D60").Value = 'Blue Cain'
$players.Range("H60").Value = 11
$players.Range("I60").Value = 8
$players.Range("J60").Value = 4
$players.Range("K60").Value = 2
$players.Range("L60").Value = 2
$players.Range("M60").Value = 1
$players.Range("N60").Value = 0
$players.Range("O60").Value = 30
# Row 63
$players.Range("D63").Value = 'Marcus Millender'
$players.Range("E63").Value = 'UGA'
$players.Range("F63").Value = 'UGA@FLA'
$players.Range("I63").Value = 18
$players.Range("J63").Value = 3
$players.Range("K63").Value = 1
$players.Range("L63").Value = 2
$players.Range("M63").Value = 0
$players.Range("N63").Value = 1
$players.Range("O63").Value = 28
# Row 64
$players.Range("D64").Value = 'Rylan Griffen'
$players.Range("E64").Value = 'TA&M'
$players.Range("F64").Value = 'TA&M@AUB'
$players.Range("H64").Value = 18
$players.Range("I64").Value = 9
$players.Range("J64").Value = 6
$players.Range("K64").Value = 1
$players.Range("L64").Value = 3
$players.Range("M64").Value = 1
$players.Range("O64").Value = 31
# Row 67
$players.Range("D67").Value = 'Mouhamed Dioubate'
$players.Range("E67").Value = 'UK'
$players.Range("F67").Value = 'MIZ@UK'
$players.Range("H67").Value = 5
$players.Range("I67").Value = 8
$players.Range("J67").Value = 2
$players.Range("K67").Value = 2
$players.Range("L67").Value = 0
$players.Range("M67").Value = 0
$players.Range("N67").Value = 3
$players.Range("O67").Value = 25
# Row 69
$players.Range("D69").Value = 'Mohamed Wague'
$players.Range("E69").Value = 'OU'
$players.Range("F69").Value = 'OU@MSST'
$players.Range("H69").Value = 7
$players.Range("I69").Value = 2
$players.Range("J69").Value = 5
$players.Range("K69").Value = 0
$players.Range("L69").Value = 0
$players.Range("M69").Value = 2
$players.Range("N69").Value = 1
$players.Range("O69").Value = 17
# Row 72
$players.Range("D72").Value = 'Mike Sharavjamts'
$players.Range("E72").Value = 'SC'
$players.Range("H72").Value = 15
$players.Range("I72").Value = 7
$players.Range("J72").Value = 7
$players.Range("K72").Value = 3
$players.Range("L72").Value = 1
$players.Range("M72").Value = 2
$players.Range("N72").Value = 0
$players.Range("O72").Value = 36
# Row 74
$players.Range("D74").Value = 'Elyjah Freeman'
$players.Range("E74").Value = 'AUB'
$players.Range("F74").Value = 'TA&M@AUB'
$players.Range("H74").Value = 10
$players.Range("I74").Value = 8
$players.Range("J74").Value = 2
$players.Range("K74").Value = 1
$players.Range("L74").Value = 1
$players.Range("M74").Value = 1
$players.Range("N74").Value = 0
$players.Range("O74").Value = 15
# Row 75
$players.Range("D75").Value = 'Tahaad Pettiford'
$players.Range("H75").Value = 8
$players.Range("I75").Value = 11
$players.Range("J75").Value = 4
$players.Range("K75").Value = 2
$players.Range("L75").Value = 0
$players.Range("M75").Value = 0
$players.Range("N75").Value = 2
$players.Range("O75").Value = 32
# Row 76
$players.Range("D76").Value = 'Max Mackinnon'
$players.Range("E76").Value = 'LSU'
$players.Range("F76").Value = 'SC@LSU'
$players.Range("H76").Value = 7
$players.Range("I76").Value = 15
$players.Range("J76").Value = 1
$players.Range("K76").Value = 0
$players.Range("L76").Value = 1
$players.Range("M76").Value = 1
$players.Range("N76").Value = 3
$players.Range("O76").Value = 27

# --- OwnerTotals sheet ---
$totals = $wb.Worksheets.Item("OwnerTotals")

# Row 2
$totals.Range("B2").Value = 130
$totals.Range("C2").Value = 5
# Row 5
$totals.Range("B5").Value = 71
$totals.Range("C5").Value = 4
# Row 6
$totals.Range("A6").Value = 'Tar'
$totals.Range("B6").Value = 69
$totals.Range("C6").Value = 4
# Row 7
$totals.Range("B7").Value = 55
$totals.Range("C7").Value = 5
# Row 8
$totals.Range("A8").Value = 'Booz'
$totals.Range("B8").Value = 51
$totals.Range("C8").Value = 5
